$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force text storage (matches source inlineStr/t="s" text cells,
# avoids Excel auto-converting numeric-looking strings to numbers),
# then reset the style so no stray number-format gets attached to the cell.

$ws.Range("D2").Value = "'42.908.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -5.05%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.216.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -6.66%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'315.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.07%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'98.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -8.88%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.585"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -7.27%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  -0.05%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.560"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -9.08%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'36.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -10.81%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'53.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.98%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.0827"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -9.94%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'7.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -8.55%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "'  -3.81%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'2.554.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -6.79%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.860"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -12.10%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'14.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -7.45%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.212.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -7.40%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'42.789.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -5.49%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'14.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.80%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'0.0₃0960"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -9.39%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'6.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -12.17%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'65.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -11.06%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'3.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -9.77%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'236.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -9.10%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'2.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -8.65%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  -0.17%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'10.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -10.13%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  -5.27%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'6.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -12.79%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'20.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -8.57%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.0877"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -9.42%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'34.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -9.57%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'156.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -7.64%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  -5.81%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'3.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +8.22%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'1.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +11.75%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.122"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -6.66%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'4.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -6.95%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  -11.96%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'3.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -5.48%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.0325"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -8.22%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'1.883.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.46%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  +0.13%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'12.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -5.17%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'88.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -10.82%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.206"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -10.12%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'5.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.75%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'60.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -13.43%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'75.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -7.67%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "'SEI"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.852"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +14.53%  "
$ws.Range("E51").Style = "Normal"
